$wb = $excel.ActiveWorkbook

# Update "展览" sheet (row5/row6 "want to go" counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 913
$ws1.Range("F6").Value = 216

# Update "全部类型" sheet (same underlying data duplicated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 913
$ws4.Range("F6").Value = 216
